# Adds the ifoCAST full-series evaluation column.
# For each data row (2-20) the existing error values shift one column to the
# left (the old "near-zero" first value is dropped), and the trailing cell
# that is no longer needed is cleared. Rows 2-10 gain one brand-new data
# point (ifoCAST) in column K that did not exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values appended at the end of the (still full-width) rows 2-10, column K.
$newK = @{
    2  = 0.3220726034063205
    3  = -1.677319868393072
    4  = 1.920156118130258
    5  = 1.480667296544459
    6  = -0.2804276406117954
    7  = 0.1301918964218456
    8  = 0.3403798556124878
    9  = -0.01153110738878937
    10 = -0.1586151035472806
}

# Last populated data column (before the edit) for each row, 2-20.
$lastCol = @{
    2  = 11  # K
    3  = 11
    4  = 11
    5  = 11
    6  = 11
    7  = 11
    8  = 11
    9  = 11
    10 = 11
    11 = 11
    12 = 10  # J
    13 = 9   # I
    14 = 8   # H
    15 = 7   # G
    16 = 6   # F
    17 = 5   # E
    18 = 4   # D
    19 = 3   # C
    20 = 2   # B
}

for ($r = 2; $r -le 20; $r++) {
    $last = $lastCol[$r]

    # Shift columns B..(last-1) to hold the value that was previously one
    # column to the right (read all old values first, then write).
    $oldVals = @{}
    for ($c = 2; $c -le $last; $c++) {
        $oldVals[$c] = $ws.Cells.Item($r, $c).Value2
    }

    for ($c = 2; $c -lt $last; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $oldVals[$c + 1]
    }

    # Rows 2-10 get a brand new ifoCAST data point in the now-vacated last
    # column instead of being cleared.
    if ($newK.ContainsKey($r)) {
        $ws.Cells.Item($r, $last).Value2 = $newK[$r]
    } else {
        $ws.Cells.Item($r, $last).ClearContents()
    }
}
